$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete row 28 first (SC 92), then row 26 (RM 232) so row indices stay valid.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()
